$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.986.83'
$ws.Range("E2").Value = '  -2.42%  '
$ws.Range("D3").Value = '2.426.29'
$ws.Range("E3").Value = '  -1.26%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '571.66'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.19%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.46'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.46%  '
$ws.Range("E7").Value = '  +0.19%  '
$ws.Range("E8").Value = '  -0.95%  '
$ws.Range("D9").Value = '2.412.65'
$ws.Range("E9").Value = '  -1.68%  '
$ws.Range("E10").Value = '  -1.05%  '
$ws.Range("E11").Value = '  -0.06%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.11'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.13%  '
$ws.Range("E13").Value = '  -2.19%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.16'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.59%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000171'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.37%  '
$ws.Range("D17").Value = '60.868.35'
$ws.Range("E17").Value = '  -2.30%  '
$ws.Range("D18").Value = '2.411.06'
$ws.Range("E18").Value = '  -1.70%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.64'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +6.57%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.68'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.96%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '323.31'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.85%  '
$ws.Range("E22").Value = '  -1.82%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.08'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.13%  '
$ws.Range("E24").Value = '  +0.16%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.88'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.77%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '65.03'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.22%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '587.89'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.76%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.46'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -8.75%  '
$ws.Range("D29").Value = '2.544.40'
$ws.Range("E29").Value = '  -1.18%  '
$ws.Range("D30").Value = '0.0₃0934'
$ws.Range("E30").Value = '  -4.96%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.94'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.53%  '
$ws.Range("E32").Value = '  -5.36%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.85'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.83%  '
$ws.Range("E34").Value = '  -1.77%  '
$ws.Range("E35").Value = '  +0.05%  '
$ws.Range("E36").Value = '  -3.00%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.66'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -6.21%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '152.45'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.32%  '
$ws.Range("E39").Value = '  -2.88%  '
$ws.Range("E40").Value = '  -0.78%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.16'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.45%  '
$ws.Range("E42").Value = '  +0.09%  '
$ws.Range("E43").Value = '  -2.77%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.19'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.29%  '
$ws.Range("E45").Value = '  -5.69%  '
$ws.Range("D46").Value = '0.0₆0290'
$ws.Range("E46").Value = '  +11.64%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '142.11'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.09%  '
$ws.Range("E48").Value = '  -3.90%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.592'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.88%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.70'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.52%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0506'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.85%  '
